# edit.ps1 -- applies the template-wording fixes described by the diff:
#   1. ".. SIGUIENTES MEDIDA Y COLINDANCIAS" -> "... MEDIDAS Y COLINDANCIAS"
#   2. the "LIBERAR DE CARGA ... AJENAS A {{SEXO_17}}" clause:
#         {{SEXO_17}} -> ÉSTA
#   3. the parallel "LIBERAR DE CARGA ... AJENAS A {{SEXO_16}}" clause:
#         {{SEXO_16}} -> ÉSTA
#   4. the "ASIMISMO, ... SE OBLIGA A CUBRIR LA PENA CONVENCIONAL ..." clause:
#         CORRESPONDAN "{{SEXO_12}} PROMITENTES {{SEXO_10}}"
#           -> CORRESPONDAN {{SEXO_7}}PROMITENTE {{SEXO_2}}"
#
# Each text edit is paired with a small bookmark touch-up so the hidden
# "_Hlk..." link bookmarks keep ending at the same semantic spot (right
# before the replaced token) that they do in the target revision, instead
# of trailing the freshly-typed replacement text.

$d = $word.ActiveDocument
$OpenCurly = [char]0x201C
$CloseCurly = [char]0x201D

# ---------------------------------------------------------------------
# Change 1: MEDIDA -> MEDIDAS
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "EL CUAL TIENEN LAS SIGUIENTES MEDIDA Y COLINDANCIAS",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "EL CUAL TIENEN LAS SIGUIENTES MEDIDAS Y COLINDANCIAS", 2)
Write-Host "Change1 (MEDIDA -> MEDIDAS):" $found

# ---------------------------------------------------------------------
# Change 2: {{SEXO_17}} -> ÉSTA  (bookmarks 50 / 51)
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "{{SEXO_17}}", $false, $false, $false, $false, $false, $true, 1, $false,
    "ÉSTA", 2)
Write-Host "Change2 (SEXO_17 -> ESTA):" $found

$bm50 = $d.Bookmarks.Item("_Hlk206455454")
$bm51 = $d.Bookmarks.Item("_Hlk207404752")
$fix = $d.Range($bm51.Start, $bm50.End)
$d.Bookmarks.Add("_Hlk207404752", $fix)

# ---------------------------------------------------------------------
# Change 3: {{SEXO_16}} -> ÉSTA  (bookmarks 56 / 57)
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "{{SEXO_16}}", $false, $false, $false, $false, $false, $true, 1, $false,
    "ÉSTA", 2)
Write-Host "Change3 (SEXO_16 -> ESTA):" $found

$bm56 = $d.Bookmarks.Item("_Hlk206455481")
$bm57 = $d.Bookmarks.Item("_Hlk207404821")
$fix = $d.Range($bm57.Start, $bm56.End)
$d.Bookmarks.Add("_Hlk207404821", $fix)

# ---------------------------------------------------------------------
# Change 4: "{{SEXO_12}} PROMITENTES {{SEXO_10}}" -> {{SEXO_7}}PROMITENTE {{SEXO_2}}"
#           (drop the opening quote, renumber the placeholders, singularize
#           PROMITENTES -> PROMITENTE) and slide bookmark 64 back to close
#           right before "CORRESPONDAN" instead of at the end of the clause.
# ---------------------------------------------------------------------
$search4 = "CORRESPONDAN " + $OpenCurly + "{{SEXO_12}} PROMITENTES {{SEXO_10}}" + $CloseCurly
$replace4 = "CORRESPONDAN {{SEXO_7}}PROMITENTE {{SEXO_2}}" + $CloseCurly
$rng = $d.Content
$found = $rng.Find.Execute(
    $search4, $false, $false, $false, $false, $false, $true, 1, $false,
    $replace4, 2)
Write-Host "Change4 (SEXO_12/10 -> SEXO_7/2):" $found

$rng = $d.Content
$found = $rng.Find.Execute(
    "ADICIONALES QUE CORRESPONDAN {{SEXO_7}}",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $rng.Start + ("ADICIONALES QUE ").Length

$bm64 = $d.Bookmarks.Item("_Hlk207405138")
$fix = $d.Range($bm64.Start, $splitPos)
$d.Bookmarks.Add("_Hlk207405138", $fix)

Write-Host "Done."
